# Insert a new record (2026/01/23, 金, 23, 201) at row 700, pushing the
# existing rows 700-741 down to 701-742 (dimension grows from D741 to D742).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 700..end down by one, opening up a blank row 700.
$ws.Rows.Item(700).Insert()

# Column A holds a literal text date ("2026/01/23"), not a real date value.
# Force the cell to Text format first so Excel doesn't auto-convert the
# "yyyy/mm/dd"-looking string into a date serial, then drop back to the
# sheet's default (unstyled) look to match the rest of the data rows.
$ws.Cells.Item(700, 1).NumberFormat = "@"
$ws.Cells.Item(700, 1).Value = "2026/01/23"
$ws.Cells.Item(700, 1).Style = "Normal"

$ws.Cells.Item(700, 2).Value = "金"
$ws.Cells.Item(700, 3).Value = 23
$ws.Cells.Item(700, 4).Value = 201
